# Weekly price update: a new week's data point was reported for
# "Terminal Hortofrutícola Agro Chillán" (Pepino dulce), inserted as the
# new row 3 (right after the existing row 2), pushing every following
# row down by one. The new row duplicates the surrounding record's
# attributes but carries its own date (44649) while keeping the same
# quality/volume/price figures as the (former) row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 3; existing rows 3-9 shift to 4-10.
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new weekly record.
$ws.Cells.Item(3, 1).Value = 7
$ws.Cells.Item(3, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(3, 3).Value = "Ñuble"
$ws.Cells.Item(3, 4).Value = 44649
$ws.Cells.Item(3, 5).Value = 16
$ws.Cells.Item(3, 6).Value = 100112043
$ws.Cells.Item(3, 7).Value = "Pepino dulce"
$ws.Cells.Item(3, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 60
$ws.Cells.Item(3, 11).Value = 15000
$ws.Cells.Item(3, 12).Value = 16000
$ws.Cells.Item(3, 13).Value = 15500
$ws.Cells.Item(3, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(3, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(3, 16).Value = 861
$ws.Cells.Item(3, 17).Value = 18
$ws.Cells.Item(3, 18).Value = "Hortaliza"
